$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D46").Value = "흉통(가슴통증) 감별진단"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/441"

$ws.Range("D52").Value = "결측값 다루기 : 왜 평균대체가 문제인가?"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2582&utm_source=rss&utm_medium=rss&utm_campaign=%25ea%25b2%25b0%25ec%25b8%25a1%25ea%25b0%2592-%25eb%258b%25a4%25eb%25a3%25a8%25ea%25b8%25b0-%25ec%2599%259c-%25ed%258f%2589%25ea%25b7%25a0%25eb%258c%2580%25ec%25b2%25b4%25ea%25b0%2580-%25eb%25ac%25b8%25ec%25a0%259c%25ec%259d%25b8%25ea%25b0%2580"

$wb.Save()
